$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.062.16'
$ws.Range("E2").Value = '  -6.30%  '

$ws.Range("D3").Value = '3.300.83'
$ws.Range("E3").Value = '  -5.04%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.59'
$ws.Range("E5").Value = '  -3.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.47'
$ws.Range("E6").Value = '  -2.48%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.304.50'
$ws.Range("E8").Value = '  -4.93%  '

$ws.Range("E9").Value = '  -2.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.37'
$ws.Range("E10").Value = '  -4.29%  '

$ws.Range("E11").Value = '  -5.68%  '

$ws.Range("E12").Value = '  -4.13%  '

$ws.Range("D13").Value = '3.859.53'
$ws.Range("E13").Value = '  -5.20%  '

$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("D15").Value = '3.294.71'

$ws.Range("E16").Value = '  -6.46%  '

$ws.Range("D17").Value = '60.274.56'
$ws.Range("E17").Value = '  -5.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '23.96'
$ws.Range("E18").Value = '  -4.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.58'
$ws.Range("E19").Value = '  -1.74%  '

$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.82'
$ws.Range("E21").Value = '  -11.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '348.29'
$ws.Range("E22").Value = '  -9.48%  '

$ws.Range("E23").Value = '  -2.80%  '

$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("D25").Value = '3.428.36'
$ws.Range("E25").Value = '  -5.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '68.57'

$ws.Range("E27").Value = '  -3.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.28'
$ws.Range("E29").Value = '  +2.71%  '

$ws.Range("E30").Value = '  +1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.80'
$ws.Range("E31").Value = '  -2.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.152'
$ws.Range("E32").Value = '  -2.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.08'
$ws.Range("E33").Value = '  -6.23%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = '3.327.40'
$ws.Range("E35").Value = '  -5.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.57'
$ws.Range("E36").Value = '  -1.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.28'
$ws.Range("E37").Value = '  +1.37%  '

$ws.Range("E38").Value = '  -0.62%  '

$ws.Range("E39").Value = '  -1.86%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '156.37'
$ws.Range("E40").Value = '  -3.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0746'
$ws.Range("E41").Value = '  -4.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.65'
$ws.Range("E43").Value = '  -2.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.28'
$ws.Range("E44").Value = '  -1.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.740'
$ws.Range("E45").Value = '  -7.28%  '

$ws.Range("E46").Value = '  +2.98%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.53'
$ws.Range("E47").Value = '  -5.33%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.48'
$ws.Range("E48").Value = '  -4.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.69'
$ws.Range("E49").Value = '  -0.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.81'
$ws.Range("E50").Value = '  +6.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.849'
$ws.Range("E51").Value = '  -5.62%  '
